$wb = $excel.ActiveWorkbook

# --- Sheet "3 V 0.3": convert E9 to a numeric value, then append row 10 ---
$ws1 = $wb.Worksheets.Item("3 V 0.3")

# E9 currently holds the bsecode as text ("543277"); store it as a number instead.
$ws1.Range("E9").Value = 543277

# Append the new data row (row 10).
$ws1.Range("A10").Value = "20/06/2024 08:44:49"
$ws1.Range("B10").Value = 1
$ws1.Range("C10").Value = "LXCHEM"
$ws1.Range("D10").Value = "Laxmi Organic Industries Ltd"
# bsecode stays textual on the new row (matches the source export), so force
# a text format before assigning or Excel will auto-coerce the digit string;
# switch the format back to General afterwards so no stray text-format style
# lingers on the cell.
$ws1.Range("E10").NumberFormat = "@"
$ws1.Range("E10").Value = "543277"
$ws1.Range("E10").Style = "Normal"
$ws1.Range("F10").Value = 3.53
$ws1.Range("G10").Value = 263.4
$ws1.Range("H10").Value = 5458317

# --- Sheet "DND 3 V 0.3": convert E5 to a numeric value, then append row 6 ---
$ws2 = $wb.Worksheets.Item("DND 3 V 0.3")

# E5 currently holds the bsecode as text ("532832"); store it as a number instead.
$ws2.Range("E5").Value = 532832

# Append the new data row (row 6).
$ws2.Range("A6").Value = "20/06/2024 08:44:49"
$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = "IBREALEST"
$ws2.Range("D6").Value = "Indiabulls Real Estate Limited"
# bsecode stays textual on the new row (matches the source export), so force
# a text format before assigning or Excel will auto-coerce the digit string;
# switch the format back to General afterwards so no stray text-format style
# lingers on the cell.
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = "532832"
$ws2.Range("E6").Style = "Normal"
$ws2.Range("F6").Value = 13.85
$ws2.Range("G6").Value = 155.89
$ws2.Range("H6").Value = 79315775
